# expense management Task (Ankur Yadav)
# Manage logical fare sheets: drop the now-unused "per-segment" columns
# (J:Y) on the all/no-logical-fare sheets, fix up the shifted
# Domestic/International flag, refresh the admin password used for the
# "nological" smoke test, and leave the UI focused back on that sheet.

$wb = $excel.ActiveWorkbook

# --- CalendarRules: just a stray click while reviewing, nothing else ---
$wsCal = $wb.Worksheets.Item("CalendarRules")
$wsCal.Activate()
$wsCal.Range("T19").Select()

# --- alllogical: remove the obsolete per-segment detail columns ---
$wsAll = $wb.Worksheets.Item("alllogical")
$wsAll.Activate()
$wsAll.Rows("1:2").Select()
$wsAll.Columns("J:Y").Delete()
$wsAll.Range("K2").Value() = "Domestic"

# --- nological: same column cleanup, plus refresh the admin password ---
$wsNo = $wb.Worksheets.Item("nological")
$wsNo.Activate()
$wsNo.Columns("J:Y").Delete()
$wsNo.Range("K2").Value() = "Domestic"
$wsNo.Range("H2").Value() = "Admin@456"

$dv = $wsNo.Range("H2").Validation
$dv.Formula1() = "Laxmi@1234,Quad@721,Admin@456"

$wsNo.Range("J14").Select()
